# Update the Cambodia MSME summary figures with more precise decimal values.
# The target cells already hold their numbers as text (shared strings), so a
# plain `.Value = "25.71"` assignment would be auto-coerced to a Double by
# Excel's normal "typed into a cell" parsing. Prefixing with a leading
# apostrophe forces the new value to stay text, but that also stamps the
# cell's style with Excel's "quote prefix" marker (a new style index). To
# keep the cell formatting identical to the original (style index 0, same as
# every other untouched data cell in these rows), we immediately paste the
# number format back in from an unaffected sibling cell in the same column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $newValue, $formatSourceAddr) {
    $ws.Range($cellAddr).Value = "'" + $newValue
    $ws.Range($formatSourceAddr).Copy() | Out-Null
    $ws.Range($cellAddr).PasteSpecial(-4122) | Out-Null
}

# Enterprises density (per 1000 people) — row 11
Set-TextValue "B11" "25.71" "B13"
Set-TextValue "C11" "0.88"  "C13"
Set-TextValue "D11" "26.59" "D13"

# Employment (% of total) — row 12
Set-TextValue "B12" "57.38" "B13"
Set-TextValue "C12" "17.26" "C13"
Set-TextValue "D12" "74.64" "D13"

# Enterprises (% of total) — row 14 (C14 "3.3" is unchanged)
Set-TextValue "B14" "96.52" "B13"
Set-TextValue "D14" "99.82" "D13"
